# Update "想去人数" (interested-count) values in column F across sheets
# 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types)
# to reflect the latest generated numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5668
$ws1.Range("F4").Value = 81
$ws1.Range("F6").Value = 950
$ws1.Range("F7").Value = 156
$ws1.Range("F8").Value = 2554
$ws1.Range("F10").Value = 169
$ws1.Range("F11").Value = 13
$ws1.Range("F12").Value = 85
$ws1.Range("F13").Value = 28
$ws1.Range("F14").Value = 2389
$ws1.Range("F15").Value = 405

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 106

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5668
$ws4.Range("F4").Value = 106
$ws4.Range("F5").Value = 81
$ws4.Range("F8").Value = 950
$ws4.Range("F9").Value = 156
$ws4.Range("F10").Value = 2554
$ws4.Range("F12").Value = 169
$ws4.Range("F13").Value = 13
$ws4.Range("F15").Value = 85
$ws4.Range("F16").Value = 28
$ws4.Range("F17").Value = 2389
$ws4.Range("F18").Value = 405
